# Apply updated dSF (column F) values for specific rows.
# These reflect a repull/recalculation of the "final delta-S" (dSF) data
# for a handful of sessions, per the commit message
# "repull data, push all data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    16 = 1
    32 = -1
    40 = 2
    44 = -1
    52 = 9
    56 = -2
    68 = -5
    69 = -2
    70 = 3
    72 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
